# Add a new "Swiss" worksheet after "Czech" by duplicating the Czech sheet
# (so it inherits the same layout/formatting) and then updating the
# market-specific cell values for Switzerland.

$wb = $excel.ActiveWorkbook
$czech = $wb.Worksheets.Item("Czech")

$czech.Copy($null, $czech)

$newSheet = $wb.Worksheets.Item($czech.Index + 1)
$newSheet.Name = "Swiss"

$newSheet.Range("B2").Value = "Switzerland Market"
$newSheet.Range("B4").Value = "NGC-3476/T2344/T2345"

# "Switzerland Market" is longer than the other market names and wraps to
# a second line in column B, so the row grows to fit (matches Excel's
# automatic row-height recalculation for wrapped text).
$newSheet.Rows(2).RowHeight = 28.8

$newSheet.Range("C12").Select()

# The previously active "Czech" tab now shows a whole-sheet selection
# (as left behind by a "Select All") instead of its old D9 selection.
$czech.Cells.Select()

$newSheet.Activate()
